$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'football youth compression pants'
$ws.Cells.Item(2, 1).Value = 'knee guard'
$ws.Cells.Item(3, 1).Value = 'knee sleeve basketball pair'
$ws.Cells.Item(4, 1).Value = 'soccer tights for boys'
$ws.Cells.Item(5, 1).Value = 'knee pads volleyball girls'
$ws.Cells.Item(6, 1).Value = 'medias basketball'
$ws.Cells.Item(7, 1).Value = 'boys yoga'
$ws.Cells.Item(8, 1).Value = 'baseball chart'
$ws.Cells.Item(9, 1).Value = 'basketball legs'
$ws.Cells.Item(10, 1).Value = 'youth knee pad'
$ws.Cells.Item(11, 1).Value = 'softball band'
$ws.Cells.Item(12, 1).Value = 'need pads volleyball'
$ws.Cells.Item(13, 1).Value = 'mens yoga pants capri'
$ws.Cells.Item(14, 1).Value = 'compression shorts knee length'
$ws.Cells.Item(15, 1).Value = 'men thigh compression pants'
$ws.Cells.Item(16, 1).Value = 'basketball leg compression'
$ws.Cells.Item(17, 1).Value = 'boys youth leggings'
$ws.Cells.Item(18, 1).Value = 'fabric softball'
$ws.Cells.Item(19, 1).Value = 'bump pad'
$ws.Cells.Item(20, 1).Value = 'men knee pads for work'
$ws.Cells.Item(21, 1).Value = 'black spandex pants men'
$ws.Cells.Item(22, 1).Value = 'youth volleyball kneepads'
$ws.Cells.Item(23, 1).Value = 'compression tight men'
$ws.Cells.Item(24, 1).Value = 'black compression tights men'
$ws.Cells.Item(25, 1).Value = 'basketball sleeve with pad'
$ws.Cells.Item(26, 1).Value = 'compression padded knee sleeve'
$ws.Cells.Item(27, 1).Value = 'compression pants black'
$ws.Cells.Item(28, 1).Value = 'men capris'
$ws.Cells.Item(29, 1).Value = 'calf pad'
$ws.Cells.Item(30, 1).Value = 'knee pads volleyball girls youth'
$ws.Cells.Item(31, 1).Value = 'men yoga pants'
$ws.Cells.Item(32, 1).Value = 'slim knee pads'
$ws.Cells.Item(33, 1).Value = 'hockey leg pads'
$ws.Cells.Item(34, 1).Value = 'soccer leggings men'
$ws.Cells.Item(35, 1).Value = 'knee pads volleyball large'
$ws.Cells.Item(36, 1).Value = 'softball fabric'
$ws.Cells.Item(37, 1).Value = 'knee pad thick'
$ws.Cells.Item(38, 1).Value = 'yoga for knees'
$ws.Cells.Item(39, 1).Value = 'knee pads for boys'
$ws.Cells.Item(40, 1).Value = 'knee sleeve honeycomb'
$ws.Cells.Item(41, 1).Value = 'black mens leggings'
$ws.Cells.Item(42, 1).Value = 'used softballs'
$ws.Cells.Item(43, 1).Value = 'youth tights for sports'
$ws.Cells.Item(44, 1).Value = 'boy compression leggings'
$ws.Cells.Item(45, 1).Value = 'cold compression pants'
$ws.Cells.Item(46, 1).Value = 'youth knee guards'
$ws.Cells.Item(47, 1).Value = 'paintball knee'
$ws.Cells.Item(48, 1).Value = 'knee pad construction'
$ws.Cells.Item(49, 1).Value = 'basketball knee sleeve youth'
$ws.Cells.Item(50, 1).Value = 'mens wrestling shorts'
$ws.Cells.Item(51, 1).Value = 'recovery compression tights men'
$ws.Cells.Item(52, 1).Value = 'womens sliding shorts softball'
$ws.Cells.Item(53, 1).Value = 'gym kneeling pad'
$ws.Cells.Item(54, 1).Value = 'gym leggings for men'
$ws.Cells.Item(55, 1).Value = 'durable pants'
$ws.Cells.Item(56, 1).Value = 'adult leggings'
$ws.Cells.Item(57, 1).Value = 'work wear knee pads'
$ws.Cells.Item(58, 1).Value = 'black lacrosse shorts'
$ws.Cells.Item(59, 1).Value = 'hex skin padding'
$ws.Cells.Item(60, 1).Value = 'sport leggings boys'
$ws.Cells.Item(61, 1).Value = 'need pads for construction'
$ws.Cells.Item(62, 1).Value = 'football girdle youth'
$ws.Cells.Item(63, 1).Value = '5 inch foam basketball'
$ws.Cells.Item(64, 1).Value = 'mens skin tight leggings'
$ws.Cells.Item(65, 1).Value = 'silicon knee pads'
$ws.Cells.Item(66, 1).Value = 'extra small baseball pants'
$ws.Cells.Item(67, 1).Value = 'football pads for men'
$ws.Cells.Item(68, 1).Value = 'mens kneepads'
$ws.Cells.Item(69, 1).Value = 'thigh protector men'
$ws.Cells.Item(70, 1).Value = 'long volleyball knee pads'
$ws.Cells.Item(71, 1).Value = 'paintball pants youth'
$ws.Cells.Item(72, 1).Value = 'gym leggings men'
$ws.Cells.Item(73, 1).Value = 'compression running capris'
$ws.Cells.Item(74, 1).Value = 'youth running pants boys'
$ws.Cells.Item(75, 1).Value = 'knee sleeve baseball'
$ws.Cells.Item(76, 1).Value = 'football pads youth'
$ws.Cells.Item(77, 1).Value = 'volleyball spandex shorts'
$ws.Cells.Item(78, 1).Value = 'work in baseball'
$ws.Cells.Item(79, 1).Value = 'thick yoga knee pad'
$ws.Cells.Item(80, 1).Value = 'running tights mens'
$ws.Cells.Item(81, 1).Value = 'basketballs 28 5'
$ws.Cells.Item(82, 1).Value = 'softball pants'
$ws.Cells.Item(83, 1).Value = 'youth volleyball spandex'
$ws.Cells.Item(84, 1).Value = 'work pants knee pads'
$ws.Cells.Item(85, 1).Value = 'size small baseball pants'
$ws.Cells.Item(86, 1).Value = 'boys sport tights'
$ws.Cells.Item(87, 1).Value = 'mens medium tall athletic pants'
$ws.Cells.Item(88, 1).Value = 'compression sleeve knee youth'
$ws.Cells.Item(89, 1).Value = 'compression knee sleeves for basketball'
$ws.Cells.Item(90, 1).Value = 'mens knee pads for work'
$ws.Cells.Item(91, 1).Value = 'adult knee pads for work'
$ws.Cells.Item(92, 1).Value = 'leggings cycling'
$ws.Cells.Item(93, 1).Value = 'male athletic tights'
$ws.Cells.Item(94, 1).Value = 'knee compression sleeve basketball'
$ws.Cells.Item(95, 1).Value = 'compressions knee'
$ws.Cells.Item(96, 1).Value = 'youth running tights'
$ws.Cells.Item(97, 1).Value = 'knee chart'
$ws.Cells.Item(98, 1).Value = 'long softball pants'
$ws.Cells.Item(99, 1).Value = 'compression pants men soccer'
$ws.Cells.Item(100, 1).Value = 'limber support'
